$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.097053050994873
$ws.Range("B1").Value = 2.099148511886597
$ws.Range("C1").Value = 9.161490440368652
$ws.Range("D1").Value = 2.417988777160645
$ws.Range("E1").Value = 1.298597574234009
